# Estadisticos Segundo Parcial Sin Ameca
# Insert a new "Ingles II" / "2ARHV" group row into each of the three
# statistics sheets (Estadisticos 1P, Estadisticos 2P, Estadisticos Final),
# pushing the previously-blank "Formación socioemocional IV" row down and
# the "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO" row down
# one more place (its data is preserved, just relocated to row 6).

$wb = $excel.ActiveWorkbook

$sheetNames = @("Estadisticos 1P", "Estadisticos 2P", "Estadisticos Final")

# Per-sheet values for the new row 4 ("Ingles II" / "2ARHV")
$row4Data = @{
    "Estadisticos 1P"     = @{ C = 30; D = 0; E = 13; F = 17; G = 56.67;              H = 6.5 }
    "Estadisticos 2P"     = @{ C = 30; D = 0; E = 0;  F = 30; G = 100;                H = 6.5 }
    "Estadisticos Final"  = @{ C = 30; D = 0; E = 0;  F = 30; G = 100;                H = 7.5 }
}

# Per-sheet values for row 6 (previously row 5, "MANTIENE EN OPERACIÓN
# CIRCUITOS DE CONTROL ELECTRÓNICO" / its existing group) after the shift
$row6Data = @{
    "Estadisticos 1P"     = @{ C = 19; D = 0; E = 9; F = 10; G = 52.63;               H = 6.2 }
    "Estadisticos 2P"     = @{ C = 19; D = 0; E = 5; F = 14; G = 73.68000000000001;   H = 6.2 }
    "Estadisticos Final"  = @{ C = 19; D = 0; E = 5; F = 14; G = 73.68000000000001;   H = 6.6 }
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Shift the old row 5 ("MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL
    # ELECTRÓNICO" with its existing stats) down to row 6 by inserting a
    # fresh blank row at position 5.
    $ws.Rows.Item(5).Insert()

    # Row 4 now holds the new "Ingles II" / "2ARHV" group with its stats.
    $ws.Range("A4").Value = "Ingles II"
    $ws.Range("B4").Value = "2ARHV"
    $ws.Range("C4").Value = $row4Data[$name].C
    $ws.Range("D4").Value = $row4Data[$name].D
    $ws.Range("E4").Value = $row4Data[$name].E
    $ws.Range("F4").Value = $row4Data[$name].F
    $ws.Range("G4").Value = $row4Data[$name].G
    $ws.Range("H4").Value = $row4Data[$name].H

    # Row 5 now holds "Formación socioemocional IV" / "4AEV" with no
    # students enrolled (all zeros, no % / average).
    $ws.Range("A5").Value = "Formación socioemocional IV"
    $ws.Range("B5").Value = "4AEV"
    $ws.Range("C5").Value = 0
    $ws.Range("D5").Value = 0
    $ws.Range("E5").Value = 0
    $ws.Range("F5").Value = 0

    # Row 6 keeps the "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL
    # ELECTRÓNICO" group (same text as before), group "4AEV", with the
    # stats that used to live on row 5.
    $ws.Range("A6").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
    $ws.Range("B6").Value = "4AEV"
    $ws.Range("C6").Value = $row6Data[$name].C
    $ws.Range("D6").Value = $row6Data[$name].D
    $ws.Range("E6").Value = $row6Data[$name].E
    $ws.Range("F6").Value = $row6Data[$name].F
    $ws.Range("G6").Value = $row6Data[$name].G
    $ws.Range("H6").Value = $row6Data[$name].H
}
